$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.203.05"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.686.88"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.50"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.08"
$ws.Range("E8").Value = "  +10.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  +4.74%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.925.99"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.691.11"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  +3.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.10"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "27.214.16"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.09"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.59"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.63"
$ws.Range("E23").Value = "  +4.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.32"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.44"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "1.550.30"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.605"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.946"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.34"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.76"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.834.94"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.792"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.79"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +5.91%  "
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("E50").Value = "  +7.30%  "
$ws.Range("E51").Value = "  +1.87%  "
